# Auto-generated script applying numeric value updates to the Lich_Profits workbook
# (scheduled market-data refresh across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC: 104 cell value update(s)
$ws.Range("L19").Value = 1563.8572
$ws.Range("N19").Value = -1913.8572
$ws.Range("M19").Value = -812
$ws.Range("H19").Value = 1323.5
$ws.Range("I19").Value = 987
$ws.Range("K19").Value = 987
$ws.Range("J19").Value = 1563.8572
$ws.Range("I31").Value = 8771.200000000001
$ws.Range("H31").Value = 8771.200000000001
$ws.Range("K31").Value = 26313.6
$ws.Range("M31").Value = -26083.6
$ws.Range("L46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("H46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("K54").Value = 50076
$ws.Range("M54").Value = -49590
$ws.Range("I54").Value = 50076
$ws.Range("H54").Value = 50076
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("H60").Value = 0
$ws.Range("H62").Value = 62503324
$ws.Range("I62").Value = 125001250
$ws.Range("K62").Value = 125001250
$ws.Range("M62").Value = -125000626
$ws.Range("M65").Value = -625003130
$ws.Range("I65").Value = 125001250
$ws.Range("H65").Value = 62503324
$ws.Range("K65").Value = 625006250
$ws.Range("M80").Value = -1582.2499
$ws.Range("H80").Value = 2339.389
$ws.Range("I80").Value = 860.0833
$ws.Range("K80").Value = 2580.2499
$ws.Range("M83").Value = -2748.7497
$ws.Range("H83").Value = 2339.389
$ws.Range("I83").Value = 860.0833
$ws.Range("K83").Value = 7740.7497
$ws.Range("I100").Value = 1767.8572
$ws.Range("H100").Value = 2287
$ws.Range("J100").Value = 3498.3333
$ws.Range("K100").Value = 1767.8572
$ws.Range("L100").Value = 3498.3333
$ws.Range("N100").Value = -4580.3333
$ws.Range("M100").Value = -1226.8572
$ws.Range("K107").Value = 31259122
$ws.Range("M107").Value = -31257202
$ws.Range("I107").Value = 31259122
$ws.Range("H107").Value = 45462816
$ws.Range("N113").Value = -10671.6665
$ws.Range("M113").Value = -3786.346
$ws.Range("H113").Value = 6500.9688
$ws.Range("J113").Value = 4163.6665
$ws.Range("I113").Value = 7040.346
$ws.Range("K113").Value = 7040.346
$ws.Range("L113").Value = 4163.6665
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 6541397.399999999
$ws.Range("M115").Value = -6539830.399999999
$ws.Range("L115").Value = 0
$ws.Range("H115").Value = 2180465.8
$ws.Range("I115").Value = 2180465.8
$ws.Range("M116").Value = -19972.143
$ws.Range("K116").Value = 23414.143
$ws.Range("I116").Value = 23414.143
$ws.Range("H116").Value = 19544.111
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("M132").Value = -9499.286
$ws.Range("I132").Value = 4009.762
$ws.Range("J132").Value = 7497.4
$ws.Range("H132").Value = 4380.787
$ws.Range("L132").Value = 22492.2
$ws.Range("K132").Value = 12029.286
$ws.Range("N132").Value = -27552.2
$ws.Range("K135").Value = 65128.2363
$ws.Range("L135").Value = 48375
$ws.Range("H135").Value = 6881.905
$ws.Range("N135").Value = -53445
$ws.Range("I135").Value = 7236.4707
$ws.Range("M135").Value = -62593.2363
$ws.Range("J135").Value = 5375
$ws.Range("K137").Value = 7149714.600000001
$ws.Range("M137").Value = -7147164.600000001
$ws.Range("H137").Value = 1697437
$ws.Range("I137").Value = 2383238.2
$ws.Range("K138").Value = 5987.2221
$ws.Range("J138").Value = 2858.5754
$ws.Range("I138").Value = 1995.7407
$ws.Range("N138").Value = -18855.7262
$ws.Range("L138").Value = 8575.726200000001
$ws.Range("H138").Value = 2625.61
$ws.Range("M138").Value = -847.2221
$ws.Range("J141").Value = 3000
$ws.Range("H141").Value = 859.75
$ws.Range("K141").Value = 2372.12898
$ws.Range("L141").Value = 9000
$ws.Range("I141").Value = 790.70966
$ws.Range("N141").Value = -19360
$ws.Range("M141").Value = 2807.87102
# ALC: 6 cell clear(s)
$ws.Range("N46").ClearContents()
$ws.Range("M46").ClearContents()
$ws.Range("M60").ClearContents()
$ws.Range("N60").ClearContents()
$ws.Range("N115").ClearContents()
$ws.Range("N131").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# ARM: 42 cell value update(s)
$ws.Range("I31").Value = 2335.5
$ws.Range("H31").Value = 4890.3335
$ws.Range("K31").Value = 2335.5
$ws.Range("M31").Value = -2041.5
$ws.Range("H32").Value = 564.49
$ws.Range("I32").Value = 560.19586
$ws.Range("K32").Value = 560.19586
$ws.Range("M32").Value = -273.19586
$ws.Range("K61").Value = 1753.8148
$ws.Range("I61").Value = 1753.8148
$ws.Range("M61").Value = -1541.8148
$ws.Range("H61").Value = 2499.919
$ws.Range("I74").Value = 1708.9803
$ws.Range("H74").Value = 2076.638
$ws.Range("K74").Value = 1708.9803
$ws.Range("M74").Value = -834.9802999999999
$ws.Range("K77").Value = 8544.9015
$ws.Range("M77").Value = -4176.9015
$ws.Range("H77").Value = 2076.638
$ws.Range("I77").Value = 1708.9803
$ws.Range("I110").Value = 1747.3334
$ws.Range("L110").Value = 4042.8572
$ws.Range("M110").Value = 297.6666
$ws.Range("N110").Value = -8132.8572
$ws.Range("J110").Value = 4042.8572
$ws.Range("H110").Value = 2983.3845
$ws.Range("K110").Value = 1747.3334
$ws.Range("I122").Value = 2353.7896
$ws.Range("H122").Value = 2353.7896
$ws.Range("K122").Value = 7061.3688
$ws.Range("M122").Value = -4611.3688
$ws.Range("M132").Value = -6156.304400000001
$ws.Range("I132").Value = 2895.4348
$ws.Range("J132").Value = 6163
$ws.Range("H132").Value = 3272.4614
$ws.Range("L132").Value = 18489
$ws.Range("K132").Value = 8686.304400000001
$ws.Range("N132").Value = -23549
$ws.Range("I136").Value = 1753.8148
$ws.Range("K136").Value = 5261.4444
$ws.Range("M136").Value = -2711.4444
$ws.Range("H136").Value = 2499.919

$ws = $wb.Worksheets.Item("BSM")
# BSM: 40 cell value update(s)
$ws.Range("H3").Value = 1617.5333
$ws.Range("M3").Value = -1261.9642
$ws.Range("K3").Value = 1375.9642
$ws.Range("I3").Value = 1375.9642
$ws.Range("K20").Value = 10846.571
$ws.Range("J20").Value = 2050.7778
$ws.Range("L20").Value = 2050.7778
$ws.Range("N20").Value = -2544.7778
$ws.Range("H20").Value = 7404.7393
$ws.Range("M20").Value = -10599.571
$ws.Range("I20").Value = 10846.571
$ws.Range("N94").Value = -2578.5555
$ws.Range("I94").Value = 1059.6
$ws.Range("M94").Value = -608.5999999999999
$ws.Range("J94").Value = 1676.5555
$ws.Range("H94").Value = 1290.9584
$ws.Range("K94").Value = 1059.6
$ws.Range("L94").Value = 1676.5555
$ws.Range("K99").Value = 3264.3704
$ws.Range("M99").Value = -1766.3704
$ws.Range("I99").Value = 3264.3704
$ws.Range("J99").Value = 3812.5
$ws.Range("N99").Value = -6808.5
$ws.Range("L99").Value = 3812.5
$ws.Range("H99").Value = 3335.0967
$ws.Range("M105").Value = -100.9231
$ws.Range("H105").Value = 1865.8572
$ws.Range("K105").Value = 1847.9231
$ws.Range("I105").Value = 1847.9231
$ws.Range("K107").Value = 533.3333
$ws.Range("M107").Value = 1386.6667
$ws.Range("I107").Value = 533.3333
$ws.Range("H107").Value = 511.36365
$ws.Range("N134").Value = -16721.3079
$ws.Range("J134").Value = 3883.7693
$ws.Range("I134").Value = 1710.9231
$ws.Range("K134").Value = 5132.7693
$ws.Range("M134").Value = -2597.7693
$ws.Range("H134").Value = 2435.205
$ws.Range("L134").Value = 11651.3079

$ws = $wb.Worksheets.Item("CRP")
# CRP: 118 cell value update(s)
$ws.Range("K7").Value = 3650
$ws.Range("M7").Value = -3537
$ws.Range("H7").Value = 7369.5
$ws.Range("I7").Value = 3650
$ws.Range("I16").Value = 1955.2222
$ws.Range("K16").Value = 1955.2222
$ws.Range("H16").Value = 1970.75
$ws.Range("M16").Value = -1668.2222
$ws.Range("K22").Value = 41804.684
$ws.Range("J22").Value = 57468.75
$ws.Range("I22").Value = 41804.684
$ws.Range("L22").Value = 57468.75
$ws.Range("H22").Value = 45981.766
$ws.Range("M22").Value = -41454.684
$ws.Range("N22").Value = -58168.75
$ws.Range("M23").Value = -902.5
$ws.Range("J23").Value = 1847.5
$ws.Range("K23").Value = 1142.5
$ws.Range("N23").Value = -2327.5
$ws.Range("L23").Value = 1847.5
$ws.Range("H23").Value = 1495
$ws.Range("I23").Value = 1142.5
$ws.Range("J27").Value = 1847.5
$ws.Range("H27").Value = 1495
$ws.Range("N27").Value = -2231.5
$ws.Range("M27").Value = -950.5
$ws.Range("K27").Value = 1142.5
$ws.Range("I27").Value = 1142.5
$ws.Range("L27").Value = 1847.5
$ws.Range("I31").Value = 1161.7826
$ws.Range("N31").Value = -3603.4614
$ws.Range("H31").Value = 1569.7797
$ws.Range("K31").Value = 1161.7826
$ws.Range("J31").Value = 3013.4614
$ws.Range("M31").Value = -866.7826
$ws.Range("L31").Value = 3013.4614
$ws.Range("H34").Value = 1569.7797
$ws.Range("K34").Value = 1161.7826
$ws.Range("M34").Value = -959.7826
$ws.Range("I34").Value = 1161.7826
$ws.Range("L34").Value = 3013.4614
$ws.Range("J34").Value = 3013.4614
$ws.Range("N34").Value = -3417.4614
$ws.Range("N53").Value = -46713.5
$ws.Range("H53").Value = 45499.5
$ws.Range("L53").Value = 45499.5
$ws.Range("J53").Value = 45499.5
$ws.Range("H54").Value = 16149.5
$ws.Range("J54").Value = 16149.5
$ws.Range("L54").Value = 16149.5
$ws.Range("N54").Value = -17465.5
$ws.Range("H58").Value = 2379.4211
$ws.Range("M58").Value = -2176.4211
$ws.Range("I58").Value = 2379.4211
$ws.Range("K58").Value = 2379.4211
$ws.Range("N74").Value = -80414
$ws.Range("J74").Value = 78666
$ws.Range("H74").Value = 78666
$ws.Range("L74").Value = 78666
$ws.Range("J77").Value = 78666
$ws.Range("H77").Value = 78666
$ws.Range("N77").Value = -244734
$ws.Range("L77").Value = 235998
$ws.Range("L97").Value = 65000
$ws.Range("N97").Value = -66982
$ws.Range("H97").Value = 65000
$ws.Range("J97").Value = 65000
$ws.Range("K99").Value = 14499.2
$ws.Range("M99").Value = -13001.2
$ws.Range("I99").Value = 14499.2
$ws.Range("J99").Value = 22501
$ws.Range("N99").Value = -25497
$ws.Range("L99").Value = 22501
$ws.Range("H99").Value = 18500.1
$ws.Range("H100").Value = 59500
$ws.Range("J100").Value = 59500
$ws.Range("L100").Value = 59500
$ws.Range("N100").Value = -61664
$ws.Range("M105").Value = -247.5
$ws.Range("H105").Value = 6315.0835
$ws.Range("K105").Value = 1994.5
$ws.Range("I105").Value = 1994.5
$ws.Range("J107").Value = 6636
$ws.Range("L107").Value = 6636
$ws.Range("H107").Value = 3950.6052
$ws.Range("N107").Value = -10476
$ws.Range("M113").Value = 214.7778000000001
$ws.Range("H113").Value = 1970.75
$ws.Range("I113").Value = 1955.2222
$ws.Range("K113").Value = 1955.2222
$ws.Range("I122").Value = 4996.6665
$ws.Range("K122").Value = 14989.9995
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("H122").Value = 4996.6665
$ws.Range("M122").Value = -12539.9995
$ws.Range("J126").Value = 22501
$ws.Range("M126").Value = -41027.60000000001
$ws.Range("H126").Value = 18500.1
$ws.Range("I126").Value = 14499.2
$ws.Range("K126").Value = 43497.60000000001
$ws.Range("N126").Value = -72443
$ws.Range("L126").Value = 67503
$ws.Range("M132").Value = -3250.7186
$ws.Range("I132").Value = 1926.9062
$ws.Range("H132").Value = 1926.9062
$ws.Range("K132").Value = 5780.7186
$ws.Range("J134").Value = 1250
$ws.Range("I134").Value = 1155.6923
$ws.Range("K134").Value = 3467.0769
$ws.Range("M134").Value = -932.0769
$ws.Range("H134").Value = 1162.4286
$ws.Range("L134").Value = 3750
$ws.Range("N134").Value = -8820
$ws.Range("I136").Value = 2379.4211
$ws.Range("K136").Value = 7138.263300000001
$ws.Range("M136").Value = -4588.263300000001
$ws.Range("H136").Value = 2379.4211
# CRP: 1 cell clear(s)
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# CUL: 89 cell value update(s)
$ws.Range("M5").Value = -1872.2
$ws.Range("H5").Value = 1367.5555
$ws.Range("K5").Value = 1984.2
$ws.Range("I5").Value = 661.4
$ws.Range("I12").Value = 606
$ws.Range("N12").Value = -16839.7861
$ws.Range("L12").Value = 16493.7861
$ws.Range("K12").Value = 1818
$ws.Range("H12").Value = 3345.48
$ws.Range("J12").Value = 5497.9287
$ws.Range("M12").Value = -1645
$ws.Range("N32").Value = -38063
$ws.Range("H32").Value = 4499.1665
$ws.Range("J32").Value = 12499
$ws.Range("I32").Value = 499.25
$ws.Range("L32").Value = 37497
$ws.Range("K32").Value = 1497.75
$ws.Range("M32").Value = -1214.75
$ws.Range("H34").Value = 968.4286
$ws.Range("K34").Value = 2837.3571
$ws.Range("M34").Value = -2753.3571
$ws.Range("I34").Value = 945.7857
$ws.Range("L34").Value = 3177
$ws.Range("J34").Value = 1059
$ws.Range("N34").Value = -3345
$ws.Range("J37").Value = 55655092
$ws.Range("N37").Value = -166965500
$ws.Range("H37").Value = 55655092
$ws.Range("L37").Value = 166965276
$ws.Range("K39").Value = 2700
$ws.Range("I39").Value = 900
$ws.Range("H39").Value = 1221.5927
$ws.Range("M39").Value = -2406
$ws.Range("N40").Value = -538
$ws.Range("L40").Value = 400
$ws.Range("J40").Value = 100
$ws.Range("M40").Value = -329
$ws.Range("H40").Value = 99.75
$ws.Range("K40").Value = 398
$ws.Range("I40").Value = 99.5
$ws.Range("H51").Value = 8749.875
$ws.Range("I51").Value = 15999.75
$ws.Range("L51").Value = 4500
$ws.Range("J51").Value = 1500
$ws.Range("M51").Value = -47539.25
$ws.Range("K51").Value = 47999.25
$ws.Range("N51").Value = -5420
$ws.Range("J52").Value = 10000
$ws.Range("N52").Value = -30532
$ws.Range("H52").Value = 10000
$ws.Range("L52").Value = 30000
$ws.Range("I57").Value = 2785
$ws.Range("K57").Value = 8355
$ws.Range("M57").Value = -7796
$ws.Range("H57").Value = 1714
$ws.Range("I87").Value = 1499
$ws.Range("H87").Value = 1499
$ws.Range("K87").Value = 4497
$ws.Range("M87").Value = -3249
$ws.Range("I90").Value = 1499
$ws.Range("K90").Value = 13491
$ws.Range("H90").Value = 1499
$ws.Range("M90").Value = -7251
$ws.Range("L97").Value = 1500
$ws.Range("N97").Value = -2492
$ws.Range("H97").Value = 317.8
$ws.Range("J97").Value = 500
$ws.Range("N113").Value = -7119.09095
$ws.Range("M113").Value = -23.60000000000036
$ws.Range("H113").Value = 865.375
$ws.Range("J113").Value = 926.36365
$ws.Range("I113").Value = 731.2
$ws.Range("K113").Value = 2193.6
$ws.Range("L113").Value = 2779.09095
$ws.Range("N129").Value = -15677.7502
$ws.Range("H129").Value = 12347261
$ws.Range("I129").Value = 22223556
$ws.Range("M129").Value = -66665668
$ws.Range("K129").Value = 66670668
$ws.Range("L129").Value = 5677.7502
$ws.Range("J129").Value = 1892.5834
$ws.Range("N134").Value = -19351.7145
$ws.Range("J134").Value = 3070.5715
$ws.Range("H134").Value = 2958.0557
$ws.Range("L134").Value = 9211.7145
$ws.Range("K135").Value = 5952.599999999999
$ws.Range("H135").Value = 1367.5555
$ws.Range("I135").Value = 661.4
$ws.Range("M135").Value = -3417.599999999999

$ws = $wb.Worksheets.Item("GSM")
# GSM: 28 cell value update(s)
$ws.Range("J39").Value = 25000
$ws.Range("H39").Value = 25000
$ws.Range("N39").Value = -26064
$ws.Range("L39").Value = 25000
$ws.Range("H70").Value = 8936.333000000001
$ws.Range("M70").Value = -8702.888999999999
$ws.Range("I70").Value = 8972.888999999999
$ws.Range("K70").Value = 8972.888999999999
$ws.Range("H73").Value = 8936.333000000001
$ws.Range("K73").Value = 8972.888999999999
$ws.Range("M73").Value = -8036.888999999999
$ws.Range("I73").Value = 8972.888999999999
$ws.Range("I102").Value = 2685.0625
$ws.Range("K102").Value = 2685.0625
$ws.Range("M102").Value = -1063.0625
$ws.Range("H102").Value = 2732.7817
$ws.Range("I122").Value = 1998.5
$ws.Range("H122").Value = 1998.5
$ws.Range("K122").Value = 5995.5
$ws.Range("M122").Value = -3545.5
$ws.Range("J126").Value = 4038
$ws.Range("H126").Value = 4678.6665
$ws.Range("N126").Value = -17054
$ws.Range("L126").Value = 12114
$ws.Range("M132").Value = -15486.938
$ws.Range("I132").Value = 6005.646
$ws.Range("H132").Value = 5931.7793
$ws.Range("K132").Value = 18016.938

$ws = $wb.Worksheets.Item("LTW")
# LTW: 44 cell value update(s)
$ws.Range("L16").Value = 738.5
$ws.Range("I16").Value = 1548.5
$ws.Range("K16").Value = 1548.5
$ws.Range("H16").Value = 1143.5
$ws.Range("N16").Value = -1078.5
$ws.Range("M16").Value = -1378.5
$ws.Range("J16").Value = 738.5
$ws.Range("K20").Value = 4750
$ws.Range("H20").Value = 4750
$ws.Range("M20").Value = -4524
$ws.Range("I20").Value = 4750
$ws.Range("J22").Value = 3663.125
$ws.Range("L22").Value = 3663.125
$ws.Range("H22").Value = 4177.5
$ws.Range("N22").Value = -4253.125
$ws.Range("J27").Value = 3663.125
$ws.Range("H27").Value = 4177.5
$ws.Range("N27").Value = -3877.125
$ws.Range("L27").Value = 3663.125
$ws.Range("M55").Value = -497.7143
$ws.Range("I55").Value = 670.7143
$ws.Range("H55").Value = 6017.7144
$ws.Range("L55").Value = 22058.715
$ws.Range("N55").Value = -22404.715
$ws.Range("J55").Value = 22058.715
$ws.Range("K55").Value = 670.7143
$ws.Range("M132").Value = -6316.4645
$ws.Range("I132").Value = 2948.8215
$ws.Range("J132").Value = 5428.467
$ws.Range("H132").Value = 4231.3965
$ws.Range("L132").Value = 16285.401
$ws.Range("K132").Value = 8846.4645
$ws.Range("N132").Value = -21345.401
$ws.Range("L136").Value = 9865.200000000001
$ws.Range("I136").Value = 1212.4584
$ws.Range("K136").Value = 3637.3752
$ws.Range("N136").Value = -14965.2
$ws.Range("M136").Value = -1087.3752
$ws.Range("H136").Value = 1570.3793
$ws.Range("J136").Value = 3288.4
$ws.Range("J141").Value = 69783.164
$ws.Range("H141").Value = 69783.164
$ws.Range("L141").Value = 69783.164
$ws.Range("N141").Value = -80143.164

$ws = $wb.Worksheets.Item("WVR")
# WVR: 52 cell value update(s)
$ws.Range("J13").Value = 14990
$ws.Range("N13").Value = -15270
$ws.Range("H13").Value = 21663
$ws.Range("L13").Value = 14990
$ws.Range("L14").Value = 5558.276
$ws.Range("N14").Value = -5894.276
$ws.Range("H14").Value = 5393.161
$ws.Range("J14").Value = 5558.276
$ws.Range("I21").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("H32").Value = 16175.333
$ws.Range("J32").Value = 20000
$ws.Range("I32").Value = 14263
$ws.Range("L32").Value = 20000
$ws.Range("K32").Value = 14263
$ws.Range("M32").Value = -13946
$ws.Range("N32").Value = -20634
$ws.Range("K35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("H35").Value = 0
$ws.Range("K107").Value = 2403.6522
$ws.Range("M107").Value = -483.6522
$ws.Range("I107").Value = 801.2174
$ws.Range("H107").Value = 855.7586
$ws.Range("L109").Value = 38996
$ws.Range("J109").Value = 38996
$ws.Range("H109").Value = 38996
$ws.Range("N109").Value = -41770
$ws.Range("I122").Value = 1383.8235
$ws.Range("H122").Value = 1503.95
$ws.Range("K122").Value = 4151.470499999999
$ws.Range("M122").Value = -1701.470499999999
$ws.Range("L125").Value = 50000
$ws.Range("J125").Value = 50000
$ws.Range("H125").Value = 50000
$ws.Range("N125").Value = -59840
$ws.Range("M126").Value = -3999.250100000001
$ws.Range("H126").Value = 2031.9445
$ws.Range("I126").Value = 2156.4167
$ws.Range("K126").Value = 6469.250100000001
$ws.Range("M132").Value = 780.8422
$ws.Range("I132").Value = 583.0526
$ws.Range("H132").Value = 1013.38776
$ws.Range("K132").Value = 1749.1578
$ws.Range("L136").Value = 9948.6666
$ws.Range("I136").Value = 4954.6045
$ws.Range("K136").Value = 14863.8135
$ws.Range("N136").Value = -15048.6666
$ws.Range("M136").Value = -12313.8135
$ws.Range("H136").Value = 4671.0386
$ws.Range("J136").Value = 3316.2222
# WVR: 2 cell clear(s)
$ws.Range("M21").ClearContents()
$ws.Range("M35").ClearContents()

